# Rename the original (only) worksheet "Seconde Sheet" -> "third", then
# duplicate it (full data + formatting) into a brand-new worksheet "emp"
# placed right after it - mirrors "learning excel read/write": the sheet
# was renamed and a duplicate ("emp") was appended.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "third"

# Copy the sheet (keeps all values, shared-string usage and number
# formatting intact) and place the copy immediately after the original.
$ws1.Copy($null, $ws1)

$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "emp"

# Keep the first sheet ("third") as the active/selected tab, matching the
# original workbook's activeTab="0".
$ws1.Activate()
